$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'313.20"
$ws.Range("E2").Value = "'-0.79%"
$ws.Range("D3").Value = "'38.06"
$ws.Range("E3").Value = "'-3.49%"
$ws.Range("D4").Value = "'5.076"
$ws.Range("E4").Value = "'-1.19%"
$ws.Range("D5").Value = "'0.07761"
$ws.Range("E5").Value = "'-5.21%"
$ws.Range("D6").Value = "'4.345"
$ws.Range("E6").Value = "'-0.53%"
$ws.Range("D7").Value = "'1.916"
$ws.Range("E7").Value = "'-3.82%"
$ws.Range("D8").Value = "'8.190"
$ws.Range("E8").Value = "'-1.59%"
$ws.Range("D9").Value = "'0.9166"
$ws.Range("E9").Value = "'-2.32%"
$ws.Range("D10").Value = "'0.1234"
$ws.Range("E10").Value = "'-5.15%"
$ws.Range("D11").Value = "'0.1891"
$ws.Range("E11").Value = "'-4.13%"
$ws.Range("D12").Value = "'0.08823"
$ws.Range("E12").Value = "'-2.34%"
$ws.Range("D13").Value = "'0.03387"
$ws.Range("E13").Value = "'-4.19%"
$ws.Range("D14").Value = "'0.09701"
$ws.Range("E14").Value = "'-0.55%"
$ws.Range("D15").Value = "'0.001373"
$ws.Range("E15").Value = "'-2.42%"
$ws.Range("D16").Value = "'0.006009"
$ws.Range("E16").Value = "'-2.12%"
$ws.Range("D17").Value = "'3.534"
$ws.Range("E17").Value = "'-2.75%"
$ws.Range("D21").Value = "'5.026"
$ws.Range("E21").Value = "'1.36%"
$ws.Range("E22").Value = "'4.11%"
$ws.Range("D23").Value = "'0.02105"
$ws.Range("E23").Value = "'5,591.97%"
$ws.Range("D24").Value = "'0.04400"
$ws.Range("E24").Value = "'0.69%"
$ws.Range("D25").Value = "'0.001213"
$ws.Range("E25").Value = "'-2.25%"
$ws.Range("D26").Value = "'0.004248"
$ws.Range("E26").Value = "'-11.03%"
$ws.Range("D27").Value = "'0.0001350"
$ws.Range("E27").Value = "'-65.29%"
$ws.Range("D39").Value = "'0.02133"
$ws.Range("E39").Value = "'-4.95%"
$ws.Range("D40").Value = "'0.04962"
$ws.Range("E40").Value = "'-4.99%"
$ws.Range("D41").Value = "'0.007773"
$ws.Range("E41").Value = "'0.36%"
$ws.Range("D42").Value = "'0.01002"
$ws.Range("E42").Value = "'-3.35%"
$ws.Range("D43").Value = "'0.1343"
$ws.Range("E43").Value = "'-4.20%"
$ws.Range("D44").Value = "'0.002061"
$ws.Range("E44").Value = "'-1.89%"
$ws.Range("D45").Value = "'0.009680"
$ws.Range("E45").Value = "'9.05%"
$ws.Range("D46").Value = "'0.00006508"
$ws.Range("E46").Value = "'-4.61%"
$ws.Range("E47").Value = "'0.01%"
$ws.Range("D48").Value = "'0.003198"
$ws.Range("E48").Value = "'6.87%"
$ws.Range("E49").Value = "'-0.10%"
$ws.Range("D50").Value = "'0.00002101"
$ws.Range("E50").Value = "'0.01%"
$ws.Range("D51").Value = "'0.0002001"
$ws.Range("E51").Value = "'0.01%"
